$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95:144 down to 96:145
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = Get-Date -Year 2023 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = 100112035
$ws.Cells.Item(95, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 30
$ws.Cells.Item(95, 11).Value = 32000
$ws.Cells.Item(95, 12).Value = 32000
$ws.Cells.Item(95, 13).Value = 32000
$ws.Cells.Item(95, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(95, 15).Value = "Región Metropolitana"
$ws.Cells.Item(95, 16).Value = 2133
$ws.Cells.Item(95, 17).Value = 15
$ws.Cells.Item(95, 18).Value = "Hortaliza"
